$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 78 (shifts existing rows 78-107 down to 79-108,
# and inherits the formatting of the row above, same as Excel's default row insert).
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new data record.
$ws.Range("A78").Value = 2
$ws.Range("B78").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44777
$ws.Range("D78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E78").Value = 4
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100107
$ws.Range("H78").Value = "Otros"
$ws.Range("I78").Value = 100107011
$ws.Range("J78").Value = "Tuna"
$ws.Range("K78").Value = "Sin especificar"
$ws.Range("L78").Value = "Primera"
$ws.Range("M78").Value = 160
$ws.Range("N78").Value = 21000
$ws.Range("O78").Value = 22000
$ws.Range("P78").Value = 21500
$ws.Range("Q78").Value = '$/caja 18 kilos'
$ws.Range("R78").Value = "Provincia de Limarí"
$ws.Range("S78").Value = 1194
$ws.Range("T78").Value = 18
